$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank columns before column A, shifting the existing
# Start Date / End Date TS / To / Cc / Bcc / Tag / Delay TS table right
# to columns E:K.
$ws.Range("A1:D1").EntireColumn.Insert()

# Populate the new leading columns (Greeting / Scope / Quote / Source)
# with sample data, row by row, left to right, so the new shared
# strings are appended to xl/sharedStrings.xml in this exact order.

# Header row
$ws.Range("A1").Value = "Greeting"
$ws.Range("B1").Value = "Scope"
$ws.Range("C1").Value = "#!@&*Quote (?!)"
$ws.Range("D1").Value = "Source"

# Row 2
$ws.Range("A2").Value = "Hello"
$ws.Range("B2").Value = "World"
$ws.Range("C2").Value = "You speak an infinite deal of nothing."
$ws.Range("D2").Value = "William Shakespeare"

# Row 3
$ws.Range("A3").Value = "Bonjour"
$ws.Range("B3").Value = "Monde"
$ws.Range("C3").Value = "These violent delights have violent ends and in their triumph die, like fire and powder which, as they kiss, consume."
$ws.Range("D3").Value = "William Shakespeare"

# Row 4
$ws.Range("A4").Value = "Ciao"
$ws.Range("B4").Value = "Mondo"
$ws.Range("C4").Value = "Conscience doth make cowards of us all."
$ws.Range("D4").Value = "William Shakespeare"

# Row 5
$ws.Range("A5").Value = "Hallo"
$ws.Range("B5").Value = "Welt"
$ws.Range("C5").Value = "When tyranny becomes law, rebellion becomes duty."
$ws.Range("D5").Value = "Thomas Jefferson"

# Row 6
$ws.Range("A6").Value = "Hola"
$ws.Range("B6").Value = "Mundo"
$ws.Range("C6").Value = "A true patriot will defend his country from its government."
$ws.Range("D6").Value = "Thomas Jefferson"
